$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Lichtwark deleted values" — refresh the B:E columns (subject group counts
# on row 1, CON/STR measurements on rows 2-3) with the updated figures.
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 172.98601574991676
$ws.Range("C2").Value = 222.92699410700871
$ws.Range("D2").Value = 173.78653927819192
$ws.Range("E2").Value = 221.42662787183104

$ws.Range("B3").Value = 171.66530788094812
$ws.Range("C3").Value = 226.6128897617138
$ws.Range("D3").Value = 176.8876067736812
$ws.Range("E3").Value = 217.0255034146733

# Match the author's saved selection (B1:E3 instead of the full B1:AY3).
$ws.Range("B1:E3").Select()
